$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.227.46"
$ws.Range("E2").Value = "  +4.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.451.11"
$ws.Range("E3").Value = "  +2.68%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.29"
$ws.Range("E5").Value = "  +2.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.19"
$ws.Range("E6").Value = "  +6.04%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +2.14%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.444.36"
$ws.Range("E8").Value = "  +2.72%  "

# Row 9
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.176"
$ws.Range("E10").Value = "  +7.80%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.645"
$ws.Range("E11").Value = "  +2.81%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.48"
$ws.Range("E12").Value = "  +3.46%  "

# Row 13
$ws.Range("E13").Value = "  +2.79%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.36"
$ws.Range("E14").Value = "  +3.56%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.007.03"
$ws.Range("E15").Value = "  +2.75%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.53"
$ws.Range("E16").Value = "  +2.03%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.453.77"
$ws.Range("E17").Value = "  +3.07%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.280.84"
$ws.Range("E18").Value = "  +4.36%  "

# Row 19
$ws.Range("E19").Value = "  +0.99%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.01"
$ws.Range("E20").Value = "  +2.69%  "

# Row 21
$ws.Range("E21").Value = "  +2.64%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "482.87"
$ws.Range("E22").Value = "  +5.36%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.08"
$ws.Range("E24").Value = "  +10.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.17"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.74"
$ws.Range("E26").Value = "  +4.37%  "

# Row 27
$ws.Range("E27").Value = "  -0.04%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.89"
$ws.Range("E28").Value = "  +1.61%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.93"
$ws.Range("E29").Value = "  +2.73%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.56"
$ws.Range("E30").Value = "  +3.70%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.98"
$ws.Range("E31").Value = "  +4.82%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "599.92"
$ws.Range("E32").Value = "  +5.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.60"
$ws.Range("E33").Value = "  +1.78%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "62.89"
$ws.Range("E34").Value = "  +2.93%  "

# Row 35
$ws.Range("E35").Value = "  +2.38%  "

# Row 36
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.11%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.147"
$ws.Range("E37").Value = "  +6.27%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.64"
$ws.Range("E38").Value = "  +0.50%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0780"
$ws.Range("E39").Value = "  +5.71%  "

# Row 40
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.387"
$ws.Range("E40").Value = "  +5.56%  "

# Row 41
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.48"
$ws.Range("E41").Value = "  +3.64%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.143.67"
$ws.Range("E42").Value = "  +2.42%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.92"
$ws.Range("E43").Value = "  +4.05%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.58"
$ws.Range("E44").Value = "  +6.02%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0425"
$ws.Range("E45").Value = "  +2.96%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").Value = "  +22.58%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.24"
$ws.Range("E47").Value = "  +2.96%  "

# Row 48
$ws.Range("E48").Value = "  +1.21%  "

# Row 49
$ws.Range("E49").Value = "  +7.19%  "

# Row 50
$ws.Range("E50").Value = "  +0.08%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.47"
$ws.Range("E51").Value = "  +3.38%  "
